$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to add in column R (year 2021), one per data row.
# Formatting for each new cell is copied from the corresponding Q cell
# (same row) so the new column matches the existing "year" columns,
# then the actual 2021 value is written on top of the copied one.
$values = @{
    4  = 2021
    5  = 99.4
    6  = 98.1
    7  = 99.319469393395053
    8  = 99.442213297634979
    9  = 99.1
    10 = 99.3
    11 = 99.799160124155549
    12 = 99.3
    13 = 99.538370126605429
    14 = 99.765563948945029
}

foreach ($row in 4..14) {
    $srcCell = $ws.Range("Q$row")
    $dstCell = $ws.Range("R$row")

    # Copy formatting (number format, font, borders, alignment, ...) from
    # the adjoining "2020" column so the new "2021" column looks the same.
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)

    # Now set the actual value for the new column.
    $dstCell.Value = $values[$row]
}

# Move the active selection, matching the saved view state in the workbook.
$null = $ws.Range("U4").Select()
